$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells (Wins, Losses, Ties) in row 1, copying the
# existing header style (bold, centered, bordered) from AB1.
$ws.Range("AB1").Copy($ws.Range("AC1"))
$ws.Range("AB1").Copy($ws.Range("AD1"))
$ws.Range("AB1").Copy($ws.Range("AE1"))

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Fill in the season record for every data row (2-53) with the team's
# Wins / Losses / Ties totals.
$ws.Range("AC2:AC53").Value = 73
$ws.Range("AD2:AD53").Value = 89
$ws.Range("AE2:AE53").Value = 0
